$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.803869605064392
$ws.Range("B1").Value = 2.31618857383728
$ws.Range("C1").Value = 2.512256383895874
$ws.Range("D1").Value = 6.231607437133789
$ws.Range("E1").Value = 0.7722772359848022
